$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for this market/category.
# It belongs chronologically before the current row 57, so insert a new
# row there and shift the existing rows 57:66 down to 58:67.
$ws.Rows("57:57").Insert()

# Fill the newly inserted row with the new observation. Columns that are
# constant across this sheet (market, region, category, quality, unit,
# origin, classification, etc.) are carried over unchanged.
$ws.Range("A57").Value = 3
$ws.Range("B57").Value = "Femacal de La Calera"
$ws.Range("C57").Value = "Coquimbo"
$ws.Range("D57").Value = 44776
$ws.Range("E57").Value = 5
$ws.Range("F57").Value = 100112035
$ws.Range("G57").Value = "Bruselas (repollito)"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 105
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 15500
$ws.Range("M57").Value = 15238
$ws.Range("N57").Value = "$/malla 15 kilos"
$ws.Range("O57").Value = "Provincia de Quillota"
$ws.Range("P57").Value = 1016
$ws.Range("Q57").Value = 15
$ws.Range("R57").Value = "Hortaliza"
